$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Games")

# New rows to append below the existing data (rows 26-41), columns A-E.
$data = @(
    @(25, 0, 8, 9, "10.01.2020"),
    @(26, 0, 0, 9, "10.01.2020"),
    @(27, 0, 0, 9, "10.01.2020"),
    @(28, 1, 2, 9, "10.01.2020"),
    @(29, 0, 0, 9, "10.01.2020"),
    @(30, 1, 3, 9, "10.01.2020"),
    @(31, 0, 2, 9, "10.01.2020"),
    @(32, 0, 0, 9, "10.01.2020"),
    @(33, 0, 0, 9, "10.01.2020"),
    @(34, 0, 0, 9, "10.01.2020"),
    @(35, 0, 0, 9, "10.01.2020"),
    @(36, 0, 0, 9, "10.01.2020"),
    @(37, 0, 0, 9, "10.01.2020"),
    @(38, 0, 0, 9, "10.01.2020"),
    @(39, 0, 0, 9, "10.01.2020"),
    @(40, 1, 7, 9, "10.01.2020")
)

$startRow = 26
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]

    # Column E holds a date written as plain text (e.g. "10.01.2020"),
    # matching the existing rows above. Force text formatting first so
    # Excel doesn't auto-convert the string into a date serial number,
    # then reset the style so no extra number-format style is left on
    # the cell.
    $cellE = $ws.Cells.Item($row, 5)
    $cellE.NumberFormat = "@"
    $cellE.Value = $vals[4]
    $cellE.Style = "Normal"
}
